# Weekly update: insert a new record row at row 288 (shifting the existing
# rows 288-356 down to 289-357) for the Brócoli / Macroferia Regional de
# Talca sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 288; this pushes rows 288..356 down to
# 289..357 (matches Excel's own Rows.Insert behaviour).
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with the new weekly record.
$ws.Cells.Item(288, 1).Value = 5
$ws.Cells.Item(288, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(288, 3).Value = "Maule"
$ws.Cells.Item(288, 4).Value = 44754
$ws.Cells.Item(288, 5).Value = 7
$ws.Cells.Item(288, 6).Value = 100112023
$ws.Cells.Item(288, 7).Value = "Brócoli"
$ws.Cells.Item(288, 8).Value = "Sin especificar"
$ws.Cells.Item(288, 9).Value = "Primera"
$ws.Cells.Item(288, 10).Value = 4000
$ws.Cells.Item(288, 11).Value = 800
$ws.Cells.Item(288, 12).Value = 800
$ws.Cells.Item(288, 13).Value = 800
$ws.Cells.Item(288, 14).Value = "`$/unidad"
$ws.Cells.Item(288, 15).Value = "Región del Maule"
$ws.Cells.Item(288, 16).Value = 800
$ws.Cells.Item(288, 17).Value = 1
$ws.Cells.Item(288, 18).Value = "Hortaliza"
